$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 256, shifting existing rows 256..342 down to 257..343
$ws.Rows(256).Insert()

# Populate the newly inserted row 256 with the new data point
$ws.Cells.Item(256, 1).Value = 3
$ws.Cells.Item(256, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(256, 3).Value = "Coquimbo"
$ws.Cells.Item(256, 4).Value = 44524
$ws.Cells.Item(256, 5).Value = 5
$ws.Cells.Item(256, 6).Value = 100112037
$ws.Cells.Item(256, 7).Value = "Cebollín"
$ws.Cells.Item(256, 8).Value = "Sin especificar"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 160
$ws.Cells.Item(256, 11).Value = 3000
$ws.Cells.Item(256, 12).Value = 3000
$ws.Cells.Item(256, 13).Value = 3000
$ws.Cells.Item(256, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(256, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(256, 16).Value = 83
$ws.Cells.Item(256, 17).Value = 36
$ws.Cells.Item(256, 18).Value = "Hortaliza"
